$wb = $excel.ActiveWorkbook

# Rename the existing sheet to "Outcomes"
$outcomes = $wb.Worksheets.Item("Sheet1")
$outcomes.Name = "Outcomes"

# Add the new "Investments" sheet after "Outcomes"
$investments = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $outcomes)
$investments.Name = "Investments"

# Populate the Investments sheet headers
$investments.Range("B3").Value = "Person"
$investments.Range("C3").Value = "Amount"
$investments.Range("D3").Value = "Date"

# Populate the Investments sheet data row
$investments.Range("B4").Value = "Peter Cohen"
$investments.Range("C4").Value = 10
$investments.Range("C4").NumberFormat = "$#,##0_);[Red]($#,##0)"
$investments.Range("D4").Value = 42170
$investments.Range("D4").NumberFormat = "d-mmm"

# Set selections on each sheet
$null = $outcomes.Range("B25").Select()
$null = $investments.Range("F7").Select()

# Make Investments the active tab
$null = $investments.Activate()
